# "second push LCC calc"
# Applies the updates to the embodied_emissions_systems workbook:
#  - PV ceiling-heating row (row 10) recalculated with new formulas/sources,
#    highlighted in red to mark the updated cells
#  - lifetime (col N) and its source (col O) updated for rows 11-14
#  - stray duplicate "per EBF" notes in M15/M16 removed
#  - selection cursor moved to O22

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 (PV / ceiling heating): new cost figures, interpolated from dp.py ---

# C10 keeps its "ceiling heating" text, just gets the red "updated" font
$ws.Range("C10").Font.Color = 255

# D10: embodied-emission value becomes a formula, shown with one decimal place
$ws.Range("D10").Formula = "=0.0425*975*30"
$ws.Range("D10").NumberFormat = "0.0"
$ws.Range("D10").Font.Color = 255

# E10 dimension label unchanged, just highlighted
$ws.Range("E10").Font.Color = 255

# F10: cost value becomes a formula as well
$ws.Range("F10").Formula = "=82*975*30"
$ws.Range("F10").Font.Color = 255

# G10: source of the emission figure
$ws.Range("G10").Value = "Aktualisierte Daten von Treeze"
$ws.Range("G10").Font.Color = 255

# K10: description of the cost figure
$ws.Range("K10").Value = "PV_cost_interpolated in dp.py"

# L10 keeps its "swissolar, BAPV" text (unchanged)

# --- Lifetime doubled (30 -> 60) and source changed to CRB for rows 11-14 ---

$ws.Range("N11").Value = 60
$ws.Range("O11").Value = "CRB"

$ws.Range("N12").Value = 60
$ws.Range("O12").Value = "CRB"

$ws.Range("N13").Value = 60
$ws.Range("O13").Value = "CRB"

$ws.Range("N14").Value = 60
$ws.Range("O14").Value = "CRB"

# --- Remove stray duplicate "per EBF" notes ---

$ws.Range("M15").ClearContents()
$ws.Range("M16").ClearContents()

# --- Move the active selection ---

$null = $ws.Range("O22").Select()
